# Update gh-pages to output generated at 456a3b4
# Updates the "want to go" (F) and "lowest price" (G) columns for a handful
# of rows on the "展览" sheet and the mirrored rows on the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (展览信息汇总表)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 65
$ws1.Range("G3").Value = 65
$ws1.Range("F5").Value = 13394
$ws1.Range("F12").Value = 13827
$ws1.Range("F13").Value = 14516
$ws1.Range("F25").Value = 117
$ws1.Range("F27").Value = 5558
$ws1.Range("F30").Value = 5355
$ws1.Range("F32").Value = 21
$ws1.Range("F33").Value = 139

# Sheet "全部类型" (all event types combined — mirrors 展览 plus more rows,
# so the row offsets differ after row 12)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 65
$ws4.Range("G3").Value = 65
$ws4.Range("F5").Value = 13394
$ws4.Range("F13").Value = 13827
$ws4.Range("F14").Value = 14516
$ws4.Range("F26").Value = 117
$ws4.Range("F28").Value = 5558
$ws4.Range("F31").Value = 5355
$ws4.Range("F33").Value = 21
$ws4.Range("F34").Value = 139
